# Auto-generated edit script: update cached market-board stat columns (H:N)
# across 8 sheets per the commit diff. Values with no replacement are cleared
# (matching the OOXML diff, which drops the <c> element entirely); values that
# newly appear are written for the first time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 649.8
$ws.Range("I12").Value = 649.6667
$ws.Range("K12").Value = 649.6667
$ws.Range("M12").Value = -479.6667
# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 3623.75
$ws.Range("J43").Value = 4331.6665
$ws.Range("L43").Value = 4331.6665
$ws.Range("N43").Value = -4469.6665
# Row 94 (Leve Item ID 19905)
$ws.Range("H94").Value = 3440.0833
$ws.Range("I94").Value = 3434.6365
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 3434.6365
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -2983.6365
$ws.Range("N94").Value = -4402
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1783.1818
$ws.Range("I132").Value = 1722
$ws.Range("J132").Value = 2395
$ws.Range("K132").Value = 5166
$ws.Range("L132").Value = 7185
$ws.Range("M132").Value = -2636
$ws.Range("N132").Value = -12245

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3342.6875
$ws.Range("I32").Value = 855.9286
$ws.Range("J32").Value = 20750
$ws.Range("K32").Value = 855.9286
$ws.Range("L32").Value = 20750
$ws.Range("M32").Value = -568.9286
$ws.Range("N32").Value = -21324
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1174.9286
$ws.Range("J61").Value = 1338.25
$ws.Range("L61").Value = 1338.25
$ws.Range("N61").Value = -1762.25
# Row 125 (Leve Item ID 34251)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1174.9286
$ws.Range("J136").Value = 1338.25
$ws.Range("L136").Value = 4014.75
$ws.Range("N136").Value = -9114.75

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1873.1666
$ws.Range("I86").Value = 1857.375
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 1857.375
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = -734.375
$ws.Range("N86").Value = -4245.5
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1873.1666
$ws.Range("I89").Value = 1857.375
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 9286.875
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = -3670.875
$ws.Range("N89").Value = -21229.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 99
$ws.Range("I7").Value = 45.5
$ws.Range("K7").Value = 45.5
$ws.Range("M7").Value = 67.5
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4414.3237
$ws.Range("I31").Value = 2478.0625
$ws.Range("J31").Value = 6135.4443
$ws.Range("K31").Value = 2478.0625
$ws.Range("L31").Value = 6135.4443
$ws.Range("M31").Value = -2183.0625
$ws.Range("N31").Value = -6725.4443
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4414.3237
$ws.Range("I34").Value = 2478.0625
$ws.Range("J34").Value = 6135.4443
$ws.Range("K34").Value = 2478.0625
$ws.Range("L34").Value = 6135.4443
$ws.Range("M34").Value = -2276.0625
$ws.Range("N34").Value = -6539.4443
# Row 52 (Leve Item ID 43237)
$ws.Range("H52").Value = 100780
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 100780
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 100780
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -101368
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 2289.9
$ws.Range("I99").Value = 1974.875
$ws.Range("K99").Value = 1974.875
$ws.Range("M99").Value = -476.875
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 2289.9
$ws.Range("I126").Value = 1974.875
$ws.Range("K126").Value = 5924.625
$ws.Range("M126").Value = -3454.625
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2515.2195
$ws.Range("I132").Value = 2422.75
$ws.Range("J132").Value = 3181
$ws.Range("K132").Value = 7268.25
$ws.Range("L132").Value = 9543
$ws.Range("M132").Value = -4738.25
$ws.Range("N132").Value = -14603
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1904.36
$ws.Range("I134").Value = 1209.4546
$ws.Range("J134").Value = 7000.3335
$ws.Range("K134").Value = 3628.3638
$ws.Range("L134").Value = 21001.0005
$ws.Range("M134").Value = -1093.3638
$ws.Range("N134").Value = -26071.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 422.73334
$ws.Range("J107").Value = 439.91666
$ws.Range("L107").Value = 1319.74998
$ws.Range("N107").Value = -5159.749980000001
# Row 128 (Leve Item ID 41814)
$ws.Range("H128").Value = 451424.28
$ws.Range("I128").Value = 451424.28
$ws.Range("K128").Value = 1354272.84
$ws.Range("M128").Value = -1349292.84

$ws = $wb.Worksheets.Item("GSM")
# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 25000
$ws.Range("I46").Value = 25000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 25000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -24844
$ws.Range("N46").ClearContents()
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2066.6667
$ws.Range("I80").Value = 1850
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 1850
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -852
$ws.Range("N80").Value = -4496
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2066.6667
$ws.Range("I83").Value = 1850
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 9250
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -4258
$ws.Range("N83").Value = -22484
# Row 109 (Leve Item ID 25691)
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37080
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 35403.28
$ws.Range("I132").Value = 44009.04
$ws.Range("K132").Value = 132027.12
$ws.Range("M132").Value = -129497.12

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2799.8
$ws.Range("I7").Value = 2999.75
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2999.75
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -2887.75
$ws.Range("N7").Value = -2224
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 683.65
$ws.Range("I55").Value = 751.2353000000001
$ws.Range("J55").Value = 300.66666
$ws.Range("K55").Value = 751.2353000000001
$ws.Range("L55").Value = 300.66666
$ws.Range("M55").Value = -578.2353000000001
$ws.Range("N55").Value = -646.66666
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2617.762
$ws.Range("I122").Value = 2293.818
$ws.Range("K122").Value = 6881.454000000001
$ws.Range("M122").Value = -4431.454000000001
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2799.8
$ws.Range("I126").Value = 2999.75
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8999.25
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -6529.25
$ws.Range("N126").Value = -10940
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 6305.875
$ws.Range("I132").Value = 4235.5
$ws.Range("J132").Value = 8376.25
$ws.Range("K132").Value = 12706.5
$ws.Range("L132").Value = 25128.75
$ws.Range("M132").Value = -10176.5
$ws.Range("N132").Value = -30188.75
# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 140 (Leve Item ID 42503)
$ws.Range("H140").Value = 134998.5
$ws.Range("J140").Value = 69997
$ws.Range("L140").Value = 69997
$ws.Range("N140").Value = -80357

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (Leve Item ID 3307)
$ws.Range("H2").Value = 9721.666999999999
$ws.Range("I2").Value = 9721.666999999999
$ws.Range("K2").Value = 9721.666999999999
$ws.Range("M2").Value = -9609.666999999999
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2039.931
$ws.Range("I122").Value = 1646.9546
$ws.Range("J122").Value = 3275
$ws.Range("K122").Value = 4940.8638
$ws.Range("L122").Value = 9825
$ws.Range("M122").Value = -2490.8638
$ws.Range("N122").Value = -14725
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2883.1667
$ws.Range("I132").Value = 2883.1667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8649.500100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6119.500100000001
$ws.Range("N132").ClearContents()
